# Update the main GSC "Video-Indexing" export (Chart sheet) with the
# latest rolling-window of daily data: drop the oldest day
# (2025-11-16, which had no data) and append three new days
# (2026-02-10, 2026-02-11, 2026-02-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the obsolete first data row (2025-11-16). This shifts every
# subsequent row up by one, so the table now starts at 2025-11-17.
$ws.Rows.Item(2).Delete()

# Helper: write a literal text value into a cell without Excel's COM
# layer re-interpreting a yyyy-mm-dd-looking string as a date serial.
# We build the text via a formula (so it's typed as text), copy it,
# and paste-special only the *values* into the target cell - this
# keeps the destination cell a plain shared-string with default
# formatting, matching how the rest of the date column is stored.
function Set-TextValue($cell, $text) {
    $scratch = $ws.Range("Z1")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.Clear()
}

# Append the new rows of daily data at the bottom of the table.
Set-TextValue $ws.Cells.Item(87, 1) "2026-02-10"
$ws.Cells.Item(87, 2).Value = 22.0
$ws.Cells.Item(87, 3).Value = 1.0
$ws.Cells.Item(87, 4).Value = 0.0

Set-TextValue $ws.Cells.Item(88, 1) "2026-02-11"
$ws.Cells.Item(88, 2).Value = 22.0
$ws.Cells.Item(88, 3).Value = 1.0
$ws.Cells.Item(88, 4).Value = 0.0

Set-TextValue $ws.Cells.Item(89, 1) "2026-02-12"
$ws.Cells.Item(89, 2).Value = 22.0
$ws.Cells.Item(89, 3).Value = 1.0
$ws.Cells.Item(89, 4).Value = 0.0
